$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.426185846328735
$ws.Range("B1").Value = 2.717191934585571
$ws.Range("C1").Value = 2.920148849487305
$ws.Range("D1").Value = 3.166586399078369
$ws.Range("E1").Value = 0.7981216311454773
